# Apply the edits described by the diff:
#  1. Dataset!C2 comment: 'int' -> 'integer'
#  2. Instructions!A10, A12, A14: point the "Standard deviation" descriptions
#     at the correct source column names instead of themselves.

$wb = $excel.ActiveWorkbook

# --- 1. Update comment text on the "Dataset" sheet, cell C2 ---
$datasetWs = $wb.Worksheets.Item("Dataset")
$null = $datasetWs.Range("C2").Comment.Text("'X' is not of type 'integer' in column 'n'")

# --- 2. Update instructional text on the "Instructions" sheet ---
# That sheet is protected, so unprotect, edit, then restore protection.
$instructionsWs = $wb.Worksheets.Item("Instructions")
$instructionsWs.Unprotect()

$instructionsWs.Range("A10").Value = "- Standard deviation in M^-1s^-1: The standard deviation of the value in 'on rate; Ka in M^-1s^-1'"
$instructionsWs.Range("A12").Value = "- Standard deviation in 1/s: The standard deviation of the value in 'off rate; Kd in 1/s'"
$instructionsWs.Range("A14").Value = "- Standard deviation in nM: The standard deviation of the value in 'dissociation constant; KD in nM'"

$instructionsWs.Protect()
